# Apply the shared-string text changes described by the diff:
#   MNF -> MFN
#   MPN -> MFP
#   WIP -> r1
#   Today -> 2020-09-09
#   8.0.8+1 -> 9.0.0
#
# These values appear identically on both the "BoM" and "DNF" worksheets:
#   D4 = Revision value, D5 = Date value, D6 = KiCad Version value
#   K8 = "MNF" column header, L8 = "MPN" column header

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "BoM" -or $ws.Name -eq "DNF") {
        $ws.Range("K8").Value = "MFN"
        $ws.Range("L8").Value = "MFP"
        $ws.Range("D4").Value = "r1"

        # "2020-09-09" looks like a date, so a plain .Value assignment would
        # get auto-converted to a date serial number (and pick up a date
        # number format) instead of staying the literal text from the diff.
        # Force the cell to Text first so it is stored as a string...
        $ws.Range("D5").NumberFormat = "@"
        $ws.Range("D5").Value = "2020-09-09"

        # ...then restore the cell's original "General" look by copying the
        # number format back from an untouched neighbor (D4), so the cell's
        # style matches the source file (only its text content changed).
        $ws.Range("D4").Copy() | Out-Null
        $ws.Range("D5").PasteSpecial(-4122) | Out-Null

        $ws.Range("D6").Value = "9.0.0"
    }
}

$excel.CutCopyMode = 0
